$d = $word.ActiveDocument

$rightApos = [char]8217   # U+2019 RIGHT SINGLE QUOTATION MARK
$enDash    = [char]8211   # U+2013 EN DASH

# The "Plan" section has a group of six bullet paragraphs, in this order:
#   0: "TF Matrix"
#   1: "Remove escape sequences "
#   2: "Identify all the unique words in all the documents (all comments) and create a TF matrix"
#   3: "TASK A: Use the TF matrix to test Ziff's law* - "
#   4: "Row sum of TF matrix will give the overall word frequency"
#   5: "Hypothesis testing for coefficient in Ziff's law equation"
# The edit effectively swaps the content of bullets 0-2 with the content of
# bullets 3-5 (plus small wording tweaks to the first two), while each
# paragraph itself (its numbering / indent level) stays in place.

$targets = @(
    "TF Matrix",
    "Remove escape sequences ",
    "Identify all the unique words in all the documents (all comments) and create a TF matrix",
    ("TASK A: Use the TF matrix to test Ziff" + $rightApos + "s law* " + $enDash + " "),
    "Row sum of TF matrix will give the overall word frequency",
    ("Hypothesis testing for coefficient in Ziff" + $rightApos + "s law equation")
)

# Find each bullet's paragraph index by matching its current text (ignoring
# the trailing paragraph mark).
$idx = @(0, 0, 0, 0, 0, 0)
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($t.Length -gt 0) {
        $t = $t.Substring(0, $t.Length - 1)
    }
    for ($k = 0; $k -lt $targets.Length; $k++) {
        if ($t -eq $targets[$k]) {
            $idx[$k] = $i
        }
    }
}

# Bullet 0 (was "TF Matrix") becomes "TASK A: Use complete corpus to test
# Ziff's law* - ", written as three runs (same formatting throughout).
$d.Paragraphs($idx[0]).Range.Text = "TASK A: Use "
$r = $d.Paragraphs($idx[0]).Range
$r.SetRange($r.End - 1, $r.End - 1)
$r.InsertAfter("complete corpus")
$r = $d.Paragraphs($idx[0]).Range
$r.SetRange($r.End - 1, $r.End - 1)
$r.InsertAfter(" to test Ziff" + $rightApos + "s law* " + $enDash + " ")

# Bullet 1 (was "Remove escape sequences ") becomes "get the overall word
# frequency", written as two runs.
$d.Paragraphs($idx[1]).Range.Text = "get"
$r = $d.Paragraphs($idx[1]).Range
$r.SetRange($r.End - 1, $r.End - 1)
$r.InsertAfter(" the overall word frequency")

# Bullet 2 (was "Identify all the unique words...") becomes "Hypothesis
# testing for coefficient in Ziff's law equation".
$d.Paragraphs($idx[2]).Range.Text = "Hypothesis testing for coefficient in Ziff" + $rightApos + "s law equation"

# Bullet 3 (was "TASK A: Use the TF matrix...") becomes "TF Matrix".
$d.Paragraphs($idx[3]).Range.Text = "TF Matrix"

# Bullet 4 (was "Row sum of TF matrix...") becomes "Remove escape sequences ".
$d.Paragraphs($idx[4]).Range.Text = "Remove escape sequences "

# Bullet 5 (was "Hypothesis testing...") becomes "Identify all the unique
# words in all the documents (all comments) and create a TF matrix".
$d.Paragraphs($idx[5]).Range.Text = "Identify all the unique words in all the documents (all comments) and create a TF matrix"

Write-Output "edit applied"
